$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as text,
# matching the source workbook where Price/Volume columns are inline strings
# (this avoids Excel auto-converting numeric-looking strings like "227.52"
# or "0.0602" into floating point numbers).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Price (column D) and Volume(1h) (column E) updates for rows whose data
# changed but whose coin identity (columns B/C) stayed the same.

# Row 2
Set-TextValue $ws.Range("D2") "37.835.79"
Set-TextValue $ws.Range("E2") "  -0.44%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.031.63"
Set-TextValue $ws.Range("E3") "  -1.15%  "

# Row 4
Set-TextValue $ws.Range("E4") "  -0.05%  "

# Row 5
Set-TextValue $ws.Range("D5") "227.52"
Set-TextValue $ws.Range("E5") "  -1.04%  "

# Row 6
Set-TextValue $ws.Range("E6") "  -0.32%  "

# Row 7
Set-TextValue $ws.Range("D7") "59.55"
Set-TextValue $ws.Range("E7") "  +2.39%  "

# Row 8
Set-TextValue $ws.Range("E8") "  -0.01%  "

# Row 9
Set-TextValue $ws.Range("E9") "  -0.35%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.0813"
Set-TextValue $ws.Range("E10") "  +0.48%  "

# Row 11
Set-TextValue $ws.Range("E11") "  +0.35%  "

# Row 12
Set-TextValue $ws.Range("D12") "14.63"
Set-TextValue $ws.Range("E12") "  +0.40%  "

# Row 13
Set-TextValue $ws.Range("E13") "  -1.22%  "

# Row 14
Set-TextValue $ws.Range("E14") "  +2.38%  "

# Row 15
Set-TextValue $ws.Range("D15") "0.762"
Set-TextValue $ws.Range("E15") "  +1.08%  "

# Row 16
Set-TextValue $ws.Range("E16") "  -1.59%  "

# Row 17
Set-TextValue $ws.Range("D17") "2.045.49"
Set-TextValue $ws.Range("E17") "  -0.34%  "

# Row 18
Set-TextValue $ws.Range("D18") "37.785.23"
Set-TextValue $ws.Range("E18") "  -0.31%  "

# Row 19
Set-TextValue $ws.Range("D19") "6.03"
Set-TextValue $ws.Range("E19") "  -1.76%  "

# Row 20
Set-TextValue $ws.Range("D20") "69.97"
Set-TextValue $ws.Range("E20") "  +0.29%  "

# Row 21
Set-TextValue $ws.Range("D21") "0.0₃0825"
Set-TextValue $ws.Range("E21") "  -0.85%  "

# Row 22
Set-TextValue $ws.Range("D22") "224.90"
Set-TextValue $ws.Range("E22") "  +0.01%  "

# Row 23
Set-TextValue $ws.Range("E23") "  +0.03%  "

# Row 24
Set-TextValue $ws.Range("D24") "2.37"
Set-TextValue $ws.Range("E24") "  -3.36%  "

# Row 25
Set-TextValue $ws.Range("D25") "2.22"
Set-TextValue $ws.Range("E25") "  -1.28%  "

# Row 26
Set-TextValue $ws.Range("E26") "  +0.47%  "

# Row 27
Set-TextValue $ws.Range("D27") "164.97"
Set-TextValue $ws.Range("E27") "  -0.86%  "

# Row 28
Set-TextValue $ws.Range("E28") "  -2.59%  "

# Row 29
Set-TextValue $ws.Range("E29") "  -0.54%  "

# Row 30
Set-TextValue $ws.Range("E30") "  -4.70%  "

# Row 31
Set-TextValue $ws.Range("D31") "0.120"
Set-TextValue $ws.Range("E31") "  +0.72%  "

# Row 32
Set-TextValue $ws.Range("D32") "4.44"
Set-TextValue $ws.Range("E32") "  -2.05%  "

# Row 33
Set-TextValue $ws.Range("D33") "2.09"
Set-TextValue $ws.Range("E33") "  +2.93%  "

# Rows 34 and 35 swapped identity (Hedera <-> InternetComputer(DFINITY)),
# along with their Link/Price/Volume data.
Set-TextValue $ws.Range("B34") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D34") "4.50"
Set-TextValue $ws.Range("E34") "  -1.70%  "

Set-TextValue $ws.Range("B35") "Hedera"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D35") "0.0602"
Set-TextValue $ws.Range("E35") "  -1.60%  "

# Row 36
Set-TextValue $ws.Range("E36") "  +6.50%  "

# Row 37
Set-TextValue $ws.Range("D37") "2.26"
Set-TextValue $ws.Range("E37") "  -3.72%  "

# Row 38
Set-TextValue $ws.Range("D38") "3.24"
Set-TextValue $ws.Range("E38") "  -2.20%  "

# Row 39
Set-TextValue $ws.Range("E39") "  +0.06%  "

# Row 40
Set-TextValue $ws.Range("D40") "1.523.69"
Set-TextValue $ws.Range("E40") "  +2.62%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.0219"
Set-TextValue $ws.Range("E41") "  +0.36%  "

# Row 42
Set-TextValue $ws.Range("D42") "96.75"
Set-TextValue $ws.Range("E42") "  -1.38%  "

# Row 43
Set-TextValue $ws.Range("D43") "16.83"
Set-TextValue $ws.Range("E43") "  +0.71%  "

# Row 45
Set-TextValue $ws.Range("E45") "  -1.73%  "

# Rows 46 and 47 swapped identity (FTXToken <-> TrustWalletToken),
# along with their Link/Price/Volume data.
Set-TextValue $ws.Range("B46") "TrustWalletToken"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D46") "1.11"
Set-TextValue $ws.Range("E46") "  -1.38%  "

Set-TextValue $ws.Range("B47") "FTXToken"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D47") "4.06"
Set-TextValue $ws.Range("E47") "  -2.75%  "

# Row 48
Set-TextValue $ws.Range("E48") "  -0.92%  "

# Row 49
Set-TextValue $ws.Range("E49") "  -0.20%  "

# Row 50
Set-TextValue $ws.Range("E50") "  +0.65%  "

# Row 51
Set-TextValue $ws.Range("D51") "2.220.10"
Set-TextValue $ws.Range("E51") "  -1.27%  "
